$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "26.309.26"
$ws.Range("E2").Value2 = "  +1.10%  "
$ws.Range("D3").Value2 = "1.610.73"
$ws.Range("E3").Value2 = "  +0.68%  "
$ws.Range("E4").Value2 = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "213.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.45%  "
$ws.Range("E6").Value2 = "  -0.15%  "
$ws.Range("E7").Value2 = "  +0.28%  "
$ws.Range("E8").Value2 = "  +0.97%  "
$ws.Range("E9").Value2 = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "18.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -0.18%  "
$ws.Range("D12").Value2 = "1.834.70"
$ws.Range("E12").Value2 = "  +0.63%  "
$ws.Range("D13").Value2 = "1.602.41"
$ws.Range("E13").Value2 = "  +0.30%  "
$ws.Range("E14").Value2 = "  +0.09%  "
$ws.Range("E15").Value2 = "  +1.18%  "
$ws.Range("D16").Value2 = "26.283.37"
$ws.Range("E16").Value2 = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "62.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +3.26%  "
$ws.Range("D18").Value2 = "0.0₃0728"
$ws.Range("E18").Value2 = "  +0.91%  "
$ws.Range("E19").Value2 = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "201.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.37%  "
$ws.Range("E21").Value2 = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "9.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +0.76%  "
$ws.Range("E23").Value2 = "  +0.62%  "
$ws.Range("E24").Value2 = "  +2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "143.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +1.53%  "
$ws.Range("E26").Value2 = "  -0.08%  "
$ws.Range("E27").Value2 = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "15.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +0.76%  "
$ws.Range("E29").Value2 = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +5.25%  "
$ws.Range("E31").Value2 = "  +0.34%  "
$ws.Range("E32").Value2 = "  +2.68%  "
$ws.Range("E33").Value2 = "  -0.26%  "
$ws.Range("E34").Value2 = "  +1.12%  "
$ws.Range("E35").Value2 = "  +1.04%  "
$ws.Range("D36").Value2 = "1.161.98"
$ws.Range("E36").Value2 = "  +3.37%  "
$ws.Range("E37").Value2 = "  +1.27%  "
$ws.Range("E38").Value2 = "  -0.08%  "
$ws.Range("E39").Value2 = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +0.42%  "
$ws.Range("E41").Value2 = "  +1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +4.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.784"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +0.19%  "
$ws.Range("D44").Value2 = "1.745.83"
$ws.Range("E44").Value2 = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "92.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -0.47%  "
$ws.Range("E46").Value2 = "  +13.07%  "
$ws.Range("E47").Value2 = "  +2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "53.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +1.15%  "
$ws.Range("E49").Value2 = "  +0.77%  "
$ws.Range("E50").Value2 = "  -0.19%  "
$ws.Range("E51").Value2 = "  -0.28%  "
